$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 9 (shifts existing rows 9.. down by one)
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with the new "enable_year" configuration entry
$ws.Cells.Item(9, 1).Value = "CHE"
$ws.Cells.Item(9, 2).Value = "trd_elecsupply"
$ws.Cells.Item(9, 3).Value = "enable_year"
$ws.Cells.Item(9, 4).Value = "configuration"
$ws.Cells.Item(9, 7).Value = 1990
$ws.Cells.Item(9, 8).Value = "yr"

# Re-apply the autofilter and named range so they cover the now-larger data range
$ws.AutoFilterMode = $false
$ws.Range("A5:L521").AutoFilter()
$wb.Names.Item("Sheet1!_FilterDatabase").RefersTo = "=Sheet1!`$A`$5:`$L`$521"

# Restore selection similar to the authored edit
$ws.Range("G23").Select()
